$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

$row = 15
$ws.Cells.Item($row, 1).Value = 30900000
$ws.Cells.Item($row, 2).Value = "EvilMage"
$ws.Cells.Item($row, 3).Value = "Normal"
$ws.Cells.Item($row, 4).Value = "Small"
$ws.Cells.Item($row, 5).Value = 20
$ws.Cells.Item($row, 6).Value = 7
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 3
$ws.Cells.Item($row, 9).Value = 3
$ws.Cells.Item($row, 10).Value = 0.3
$ws.Cells.Item($row, 11).Value = 2

$ws.Range("M16").Select() | Out-Null
